$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026692233669576
$ws.Cells.Item(2, 4).Value = 1.030942046220045
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.025213079140802
$ws.Cells.Item(2, 9).Value = 1.03404718741774
$ws.Cells.Item(2, 10).Value = 1.031854259123576
$ws.Cells.Item(2, 11).Value = 1.033751598896578
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.028039309087408
$ws.Cells.Item(2, 14).Value = 1.033319610526876

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027558117530428
$ws.Cells.Item(3, 4).Value = 1.031576227721089
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.026712086702378
$ws.Cells.Item(3, 9).Value = 1.034258300101165
$ws.Cells.Item(3, 10).Value = 1.032360441212386
$ws.Cells.Item(3, 11).Value = 1.034194708397107
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.029343694225454
$ws.Cells.Item(3, 14).Value = 1.033826511452306

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028118130816821
$ws.Cells.Item(4, 4).Value = 1.031986280638885
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.027681797559175
$ws.Cells.Item(4, 9).Value = 1.034393345810776
$ws.Cells.Item(4, 10).Value = 1.032687072150172
$ws.Cells.Item(4, 11).Value = 1.034480438178519
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.030186981760471
$ws.Cells.Item(4, 14).Value = 1.034153606243489

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028353495237646
$ws.Cells.Item(5, 4).Value = 1.032158593198752
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.028089409353528
$ws.Cells.Item(5, 9).Value = 1.034449745773814
$ws.Cells.Item(5, 10).Value = 1.032824171595647
$ws.Cells.Item(5, 11).Value = 1.03460032128279
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.030541327282527
$ws.Cells.Item(5, 14).Value = 1.034290900385902

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028393010133092
$ws.Cells.Item(6, 4).Value = 1.032187520894553
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.028157846094656
$ws.Cells.Item(6, 9).Value = 1.034459193682437
$ws.Cells.Item(6, 10).Value = 1.032847178515028
$ws.Cells.Item(6, 11).Value = 1.034620436236588
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.030600813447744
$ws.Cells.Item(6, 14).Value = 1.034313939977747

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.028121276024015
$ws.Cells.Item(7, 4).Value = 1.031988583378166
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.027687244297005
$ws.Cells.Item(7, 9).Value = 1.034394100896719
$ws.Cells.Item(7, 10).Value = 1.032688904929319
$ws.Cells.Item(7, 11).Value = 1.034482040996777
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.030191717215934
$ws.Cells.Item(7, 14).Value = 1.034155441625392

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.026984919282208
$ws.Cells.Item(8, 4).Value = 1.031156433552329
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.025719729703804
$ws.Cells.Item(8, 9).Value = 1.034118856483683
$ws.Cells.Item(8, 10).Value = 1.032025512594444
$ws.Cells.Item(8, 11).Value = 1.033901554996677
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.028480287046571
$ws.Cells.Item(8, 14).Value = 1.033491107197309

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.024980441183017
$ws.Cells.Item(9, 4).Value = 1.02968777360115
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.022250619913363
$ws.Cells.Item(9, 9).Value = 1.033621913811165
$ws.Cells.Item(9, 10).Value = 1.030849619916807
$ws.Cells.Item(9, 11).Value = 1.032871081107615
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.025458712099425
$ws.Cells.Item(9, 14).Value = 1.032313544617193

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023642735050104
$ws.Cells.Item(10, 4).Value = 1.028707154704019
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.019936182879611
$ws.Cells.Item(10, 9).Value = 1.033282609073022
$ws.Cells.Item(10, 10).Value = 1.03006104691607
$ws.Cells.Item(10, 11).Value = 1.032179013727712
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.023440174803127
$ws.Cells.Item(10, 14).Value = 1.031523851752345

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.023063164647392
$ws.Cells.Item(11, 4).Value = 1.028282184534859
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.018933539194082
$ws.Cells.Item(11, 9).Value = 1.033133788090838
$ws.Cells.Item(11, 10).Value = 1.029718484047478
$ws.Cells.Item(11, 11).Value = 1.031878137120331
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.022565085848699
$ws.Cells.Item(11, 14).Value = 1.031180802405186

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022847835985572
$ws.Cells.Item(12, 4).Value = 1.028124278766424
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.018561035373784
$ws.Cells.Item(12, 9).Value = 1.033078223980558
$ws.Cells.Item(12, 10).Value = 1.02959107471594
$ws.Cells.Item(12, 11).Value = 1.031766196932012
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.02223987644732
$ws.Cells.Item(12, 14).Value = 1.031053212137783

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022894026997386
$ws.Cells.Item(13, 4).Value = 1.028158152476606
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.018640942296839
$ws.Cells.Item(13, 9).Value = 1.033090155583708
$ws.Cells.Item(13, 10).Value = 1.029618411974504
$ws.Cells.Item(13, 11).Value = 1.031790216681485
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.022309642394163
$ws.Cells.Item(13, 14).Value = 1.03108058821839

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.023045366542003
$ws.Cells.Item(14, 4).Value = 1.028269133071504
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.01890274952186
$ws.Cells.Item(14, 9).Value = 1.033129200965357
$ws.Cells.Item(14, 10).Value = 1.029707955744214
$ws.Cells.Item(14, 11).Value = 1.031868887805736
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.022538207278905
$ws.Cells.Item(14, 14).Value = 1.031170259150524

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023138605147295
$ws.Cells.Item(15, 4).Value = 1.028337504880919
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.019064047180615
$ws.Cells.Item(15, 9).Value = 1.033153220300641
$ws.Cells.Item(15, 10).Value = 1.029763104584926
$ws.Cells.Item(15, 11).Value = 1.031917335671272
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.022679011984331
$ws.Cells.Item(15, 14).Value = 1.031225486308916

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023681192165845
$ws.Cells.Item(16, 4).Value = 1.028735351124356
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.020002714505992
$ws.Cells.Item(16, 9).Value = 1.033292445809505
$ws.Cells.Item(16, 10).Value = 1.030083758381564
$ws.Cells.Item(16, 11).Value = 1.032198956498523
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.023498229064691
$ws.Cells.Item(16, 14).Value = 1.031546595470724

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024021452955038
$ws.Cells.Item(17, 4).Value = 1.028984814694942
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.020591384006291
$ws.Cells.Item(17, 9).Value = 1.033379269702269
$ws.Cells.Item(17, 10).Value = 1.030284599957418
$ws.Cells.Item(17, 11).Value = 1.032375286763008
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.02401181807073
$ws.Cells.Item(17, 14).Value = 1.031747722264653

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024219889019002
$ws.Cells.Item(18, 4).Value = 1.02913028821967
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.020934699304118
$ws.Cells.Item(18, 9).Value = 1.033429729304476
$ws.Cells.Item(18, 10).Value = 1.030401640789912
$ws.Cells.Item(18, 11).Value = 1.032478020771928
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.024311284960166
$ws.Cells.Item(18, 14).Value = 1.031864929308554

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.024287545113615
$ws.Cells.Item(19, 4).Value = 1.029179885089075
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.021051753374173
$ws.Cells.Item(19, 9).Value = 1.033446903631425
$ws.Cells.Item(19, 10).Value = 1.030441530601575
$ws.Cells.Item(19, 11).Value = 1.032513030658857
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.02441337846386
$ws.Cells.Item(19, 14).Value = 1.031904875768325

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02398494953386
$ws.Cells.Item(20, 4).Value = 1.028958053164423
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.020528230103055
$ws.Cells.Item(20, 9).Value = 1.033369973279978
$ws.Cells.Item(20, 10).Value = 1.030263062590328
$ws.Cells.Item(20, 11).Value = 1.032356380238673
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.023956725300444
$ws.Cells.Item(20, 14).Value = 1.031726154312031

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.023000802210577
$ws.Cells.Item(21, 4).Value = 1.028236453509343
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.018825655998953
$ws.Cells.Item(21, 9).Value = 1.033117710950897
$ws.Cells.Item(21, 10).Value = 1.029681591925705
$ws.Cells.Item(21, 11).Value = 1.031845726122276
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.022470905138898
$ws.Cells.Item(21, 14).Value = 1.031143857892369

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022381738018816
$ws.Cells.Item(22, 4).Value = 1.027782449201242
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.017754729974748
$ws.Cells.Item(22, 9).Value = 1.032957452564394
$ws.Cells.Item(22, 10).Value = 1.029315036416127
$ws.Cells.Item(22, 11).Value = 1.031523609172294
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.021535769039107
$ws.Cells.Item(22, 14).Value = 1.030776781831922

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022709943348324
$ws.Cells.Item(23, 4).Value = 1.028023154307434
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.018322492672622
$ws.Cells.Item(23, 9).Value = 1.033042565051714
$ws.Cells.Item(23, 10).Value = 1.029509445643441
$ws.Cells.Item(23, 11).Value = 1.031694468811495
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.022031593278839
$ws.Cells.Item(23, 14).Value = 1.030971467142639

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024001443958249
$ws.Cells.Item(24, 4).Value = 1.028970145654249
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.020556766773256
$ws.Cells.Item(24, 9).Value = 1.033374174499368
$ws.Cells.Item(24, 10).Value = 1.030272794728877
$ws.Cells.Item(24, 11).Value = 1.032364923642923
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.023981619662461
$ws.Cells.Item(24, 14).Value = 1.031735900271333

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025498891757682
$ws.Cells.Item(25, 4).Value = 1.030067726684354
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.023147745996588
$ws.Cells.Item(25, 9).Value = 1.033751797819524
$ws.Cells.Item(25, 10).Value = 1.030849619916807
$ws.Cells.Item(25, 11).Value = 1.033138380520923
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.026240575916266
$ws.Cells.Item(25, 14).Value = 1.032618792664505
